# 16.7.1.1.xlsx — add a "2020" data column (column L) that mirrors the
# existing 2019 column (K), carrying over its style for every row, and
# select the new column's data range.
#
# Column L values (per the target diff):
#   L4  = 2020            (header year, style copied from K4)
#   L5  = 1.2              L6  = 1.7              L7  = 0.4
#   L8  = 3.3              L9  = 3.9              L10 = 2.4
#   L11 = 95.5             L12 = 94.4             L13 = 97.2
# Each L cell reuses the exact style of its K counterpart on the same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the view so column C becomes the left-most visible column
# (mirrors the author's topLeftCell="C1" view state).
$excel.ActiveWindow.ScrollColumn = 3

# --- Row 4 (header year) ---------------------------------------------
$ws.Range("K4").Copy($ws.Range("L4"))
$ws.Range("L4").Value = 2020

# --- Row 5 --------------------------------------------------------------
# L5 already existed (empty, style 16) and becomes style 15 with 1.2,
# identical to K5 -> just clone K5 wholesale onto L5.
$ws.Range("K5").Copy($ws.Range("L5"))

# --- Row 6 (new cell, same as K6) ---------------------------------------
$ws.Range("K6").Copy($ws.Range("L6"))

# --- Row 7 (new cell, style like K7 but a different value) -------------
$ws.Range("K7").Copy($ws.Range("L7"))
$ws.Range("L7").Value = 0.4

# --- Row 8 (existing empty cell -> style/value like K8) ----------------
$ws.Range("K8").Copy($ws.Range("L8"))

# --- Row 9 (new cell, same as K9) ---------------------------------------
$ws.Range("K9").Copy($ws.Range("L9"))

# --- Row 10 (new cell, style like K10 but a different value) -----------
$ws.Range("K10").Copy($ws.Range("L10"))
$ws.Range("L10").Value = 2.4

# --- Row 11 (existing empty cell -> style/value like K11) --------------
$ws.Range("K11").Copy($ws.Range("L11"))

# --- Row 12 (new cell, same as K12) -------------------------------------
$ws.Range("K12").Copy($ws.Range("L12"))

# --- Row 13 (new cell, style like K13 but a different value) -----------
$ws.Range("K13").Copy($ws.Range("L13"))
$ws.Range("L13").Value = 97.2

# Match the author's final selection: L4:L13 with the active cell at L4.
$ws.Range("L4:L13").Select()
